$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Widen column F (DSP) from 10.46875 to 11.71875 chars ---
# The COM ColumnWidth setter here quantizes the stored raw width to
# multiples of 1/6 (output = round(input to 1/6) + 5/6), so feed it a
# pre-compensated value (target - 5/6, snapped to a 1/6 step) to land on
# the closest representable raw width to 11.71875.
$ws.Columns.Item(6).ColumnWidth = 10.833333333333334

# --- Update the 64-bit-design utilization values on row 2 ---
$ws.Range("B2").Value = 22.526315689086914
$ws.Range("C2").Value = 5.344827651977539
$ws.Range("D2").Value = 14.88063907623291
$ws.Range("E2").Value = 57.85714340209961
$ws.Range("F2").Value = 57.272727966308594
